$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.822043180465698
$ws.Range("B1").Value = 4.178537845611572
$ws.Range("C1").Value = 3.191492080688477
$ws.Range("D1").Value = 2.948970317840576
$ws.Range("E1").Value = 2.646575689315796
